$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 901.60785
$ws_ALC.Range("J17").Value = 909.9
$ws_ALC.Range("L17").Value = 2729.7
$ws_ALC.Range("N17").Value = -3065.7

# ALC row 18
$ws_ALC.Range("H18").Value = 1800.1428
$ws_ALC.Range("I18").Value = 1800.1428
$ws_ALC.Range("J18").Value = 0
$ws_ALC.Range("K18").Value = 1800.1428
$ws_ALC.Range("L18").Value = 0
$ws_ALC.Range("M18").Value = -1516.1428
$ws_ALC.Range("N18").ClearContents()

# ALC row 32
$ws_ALC.Range("H32").Value = 3868.6667
$ws_ALC.Range("I32").Value = 3792.6667
$ws_ALC.Range("K32").Value = 3792.6667
$ws_ALC.Range("M32").Value = -3466.6667

# ALC row 33
$ws_ALC.Range("H33").Value = 4131.923
$ws_ALC.Range("I33").Value = 4899.524
$ws_ALC.Range("K33").Value = 4899.524
$ws_ALC.Range("M33").Value = -4670.524

# ALC row 70
$ws_ALC.Range("H70").Value = 2809.3
$ws_ALC.Range("J70").Value = 2809.3
$ws_ALC.Range("L70").Value = 8427.900000000001
$ws_ALC.Range("N70").Value = -8967.900000000001

# ALC row 73
$ws_ALC.Range("H73").Value = 2809.3
$ws_ALC.Range("J73").Value = 2809.3
$ws_ALC.Range("L73").Value = 8427.900000000001
$ws_ALC.Range("N73").Value = -10299.9

# ALC row 88
$ws_ALC.Range("H88").Value = 13526.583
$ws_ALC.Range("J88").Value = 11968.889
$ws_ALC.Range("L88").Value = 11968.889
$ws_ALC.Range("N88").Value = -12780.889

# ALC row 91
$ws_ALC.Range("H91").Value = 13526.583
$ws_ALC.Range("J91").Value = 11968.889
$ws_ALC.Range("L91").Value = 11968.889
$ws_ALC.Range("N91").Value = -14776.889

# ALC row 96
$ws_ALC.Range("H96").Value = 497
$ws_ALC.Range("I96").Value = 118.666664
$ws_ALC.Range("K96").Value = 355.999992
$ws_ALC.Range("M96").Value = 1017.000008

# ALC row 99
$ws_ALC.Range("H99").Value = 4304
$ws_ALC.Range("I99").Value = 174
$ws_ALC.Range("J99").Value = 10499
$ws_ALC.Range("K99").Value = 522
$ws_ALC.Range("L99").Value = 31497
$ws_ALC.Range("M99").Value = 976
$ws_ALC.Range("N99").Value = -34493

# ALC row 107
$ws_ALC.Range("H107").Value = 813.6
$ws_ALC.Range("I107").Value = 616.1667
$ws_ALC.Range("J107").Value = 1603.3334
$ws_ALC.Range("K107").Value = 616.1667
$ws_ALC.Range("L107").Value = 1603.3334
$ws_ALC.Range("M107").Value = 1303.8333
$ws_ALC.Range("N107").Value = -5443.3334

# ALC row 132
$ws_ALC.Range("H132").Value = 8823.299999999999
$ws_ALC.Range("I132").Value = 9291.471
$ws_ALC.Range("J132").Value = 6170.3335
$ws_ALC.Range("K132").Value = 27874.413
$ws_ALC.Range("L132").Value = 18511.0005
$ws_ALC.Range("M132").Value = -25344.413
$ws_ALC.Range("N132").Value = -23571.0005

# ALC row 137
$ws_ALC.Range("H137").Value = 30792.057
$ws_ALC.Range("I137").Value = 38760.703
$ws_ALC.Range("K137").Value = 116282.109
$ws_ALC.Range("M137").Value = -113732.109

# ALC row 141
$ws_ALC.Range("H141").Value = 1597.1428
$ws_ALC.Range("I141").Value = 1556.6
$ws_ALC.Range("K141").Value = 4669.799999999999
$ws_ALC.Range("M141").Value = 510.2000000000007

# ARM row 43
$ws_ARM.Range("H43").Value = 21591.25
$ws_ARM.Range("J43").Value = 21591.25
$ws_ARM.Range("L43").Value = 21591.25
$ws_ARM.Range("N43").Value = -22217.25

# ARM row 61
$ws_ARM.Range("H61").Value = 6332.647
$ws_ARM.Range("I61").Value = 5546.857
$ws_ARM.Range("J61").Value = 9999.666999999999
$ws_ARM.Range("K61").Value = 5546.857
$ws_ARM.Range("L61").Value = 9999.666999999999
$ws_ARM.Range("M61").Value = -5334.857
$ws_ARM.Range("N61").Value = -10423.667

# ARM row 102
$ws_ARM.Range("H102").Value = 4336.3335
$ws_ARM.Range("J102").Value = 4336.3335
$ws_ARM.Range("L102").Value = 4336.3335
$ws_ARM.Range("N102").Value = -7580.3335

# ARM row 110
$ws_ARM.Range("H110").Value = 3872.4546
$ws_ARM.Range("I110").Value = 2691.3333
$ws_ARM.Range("J110").Value = 5289.8
$ws_ARM.Range("K110").Value = 2691.3333
$ws_ARM.Range("L110").Value = 5289.8
$ws_ARM.Range("M110").Value = -646.3332999999998
$ws_ARM.Range("N110").Value = -9379.799999999999

# ARM row 122
$ws_ARM.Range("H122").Value = 1348.28
$ws_ARM.Range("I122").Value = 1301.625
$ws_ARM.Range("K122").Value = 3904.875
$ws_ARM.Range("M122").Value = -1454.875

# ARM row 132
$ws_ARM.Range("H132").Value = 31154.055
$ws_ARM.Range("I132").Value = 34110.844
$ws_ARM.Range("J132").Value = 7499.75
$ws_ARM.Range("K132").Value = 102332.532
$ws_ARM.Range("L132").Value = 22499.25
$ws_ARM.Range("M132").Value = -99802.53199999999
$ws_ARM.Range("N132").Value = -27559.25

# ARM row 136
$ws_ARM.Range("H136").Value = 6332.647
$ws_ARM.Range("I136").Value = 5546.857
$ws_ARM.Range("J136").Value = 9999.666999999999
$ws_ARM.Range("K136").Value = 16640.571
$ws_ARM.Range("L136").Value = 29999.001
$ws_ARM.Range("M136").Value = -14090.571
$ws_ARM.Range("N136").Value = -35099.001

# BSM row 60
$ws_BSM.Range("H60").Value = 69894.5
$ws_BSM.Range("J60").Value = 69894.5
$ws_BSM.Range("L60").Value = 69894.5
$ws_BSM.Range("N60").Value = -71092.5

# BSM row 134
$ws_BSM.Range("H134").Value = 2416.5435
$ws_BSM.Range("I134").Value = 2425.8
$ws_BSM.Range("J134").Value = 2000
$ws_BSM.Range("K134").Value = 7277.400000000001
$ws_BSM.Range("L134").Value = 6000
$ws_BSM.Range("M134").Value = -4742.400000000001
$ws_BSM.Range("N134").Value = -11070

# CRP row 58
$ws_CRP.Range("H58").Value = 336335
$ws_CRP.Range("I58").Value = 336335
$ws_CRP.Range("J58").Value = 0
$ws_CRP.Range("K58").Value = 336335
$ws_CRP.Range("L58").Value = 0
$ws_CRP.Range("M58").Value = -336132
$ws_CRP.Range("N58").ClearContents()

# CRP row 132
$ws_CRP.Range("H132").Value = 2748.6333
$ws_CRP.Range("I132").Value = 2591.3333
$ws_CRP.Range("K132").Value = 7773.999899999999
$ws_CRP.Range("M132").Value = -5243.999899999999

# CRP row 134
$ws_CRP.Range("H134").Value = 47799.727
$ws_CRP.Range("I134").Value = 64282
$ws_CRP.Range("K134").Value = 192846
$ws_CRP.Range("M134").Value = -190311

# CRP row 136
$ws_CRP.Range("H136").Value = 336335
$ws_CRP.Range("I136").Value = 336335
$ws_CRP.Range("J136").Value = 0
$ws_CRP.Range("K136").Value = 1009005
$ws_CRP.Range("L136").Value = 0
$ws_CRP.Range("M136").Value = -1006455
$ws_CRP.Range("N136").ClearContents()

# CUL row 11
$ws_CUL.Range("H11").Value = 355.75
$ws_CUL.Range("I11").Value = 426
$ws_CUL.Range("J11").Value = 145
$ws_CUL.Range("K11").Value = 1278
$ws_CUL.Range("L11").Value = 435
$ws_CUL.Range("M11").Value = -1138
$ws_CUL.Range("N11").Value = -715

# GSM row 80
$ws_GSM.Range("H80").Value = 3706.2666
$ws_GSM.Range("I80").Value = 3304.5
$ws_GSM.Range("J80").Value = 4165.4287
$ws_GSM.Range("K80").Value = 3304.5
$ws_GSM.Range("L80").Value = 4165.4287
$ws_GSM.Range("M80").Value = -2306.5
$ws_GSM.Range("N80").Value = -6161.4287

# GSM row 83
$ws_GSM.Range("H83").Value = 3706.2666
$ws_GSM.Range("I83").Value = 3304.5
$ws_GSM.Range("J83").Value = 4165.4287
$ws_GSM.Range("K83").Value = 16522.5
$ws_GSM.Range("L83").Value = 20827.1435
$ws_GSM.Range("M83").Value = -11530.5
$ws_GSM.Range("N83").Value = -30811.1435

# GSM row 102
$ws_GSM.Range("H102").Value = 3805.45
$ws_GSM.Range("I102").Value = 3194.375
$ws_GSM.Range("K102").Value = 3194.375
$ws_GSM.Range("M102").Value = -1572.375

# LTW row 16
$ws_LTW.Range("H16").Value = 3217.1667
$ws_LTW.Range("I16").Value = 3244.4614
$ws_LTW.Range("J16").Value = 3039.75
$ws_LTW.Range("K16").Value = 3244.4614
$ws_LTW.Range("L16").Value = 3039.75
$ws_LTW.Range("M16").Value = -3074.4614
$ws_LTW.Range("N16").Value = -3379.75

# LTW row 40
$ws_LTW.Range("H40").Value = 7222.759
$ws_LTW.Range("I40").Value = 5452.727
$ws_LTW.Range("K40").Value = 5452.727
$ws_LTW.Range("M40").Value = -5316.727

# LTW row 55
$ws_LTW.Range("H55").Value = 649.5454999999999
$ws_LTW.Range("I55").Value = 263.26666
$ws_LTW.Range("K55").Value = 263.26666
$ws_LTW.Range("M55").Value = -90.26666

# LTW row 61
$ws_LTW.Range("H61").Value = 5008.6665
$ws_LTW.Range("I61").Value = 1376.8889
$ws_LTW.Range("K61").Value = 1376.8889
$ws_LTW.Range("M61").Value = -1174.8889

# LTW row 68
$ws_LTW.Range("H68").Value = 8600

# LTW row 71
$ws_LTW.Range("H71").Value = 8600

# LTW row 113
$ws_LTW.Range("H113").Value = 5008.6665
$ws_LTW.Range("I113").Value = 1376.8889
$ws_LTW.Range("K113").Value = 1376.8889
$ws_LTW.Range("M113").Value = 793.1111000000001

# LTW row 132
$ws_LTW.Range("H132").Value = 32821.098
$ws_LTW.Range("I132").Value = 40761.72
$ws_LTW.Range("J132").Value = 4587.778
$ws_LTW.Range("K132").Value = 122285.16
$ws_LTW.Range("L132").Value = 13763.334
$ws_LTW.Range("M132").Value = -119755.16
$ws_LTW.Range("N132").Value = -18823.334

# LTW row 136
$ws_LTW.Range("H136").Value = 5786.909
$ws_LTW.Range("I136").Value = 0
$ws_LTW.Range("J136").Value = 5786.909
$ws_LTW.Range("K136").Value = 0
$ws_LTW.Range("L136").Value = 17360.727
$ws_LTW.Range("M136").ClearContents()
$ws_LTW.Range("N136").Value = -22460.727

# WVR row 62
$ws_WVR.Range("H62").Value = 172748
$ws_WVR.Range("I62").Value = 4997.5
$ws_WVR.Range("J62").Value = 256623.25
$ws_WVR.Range("K62").Value = 4997.5
$ws_WVR.Range("L62").Value = 256623.25
$ws_WVR.Range("M62").Value = -4373.5
$ws_WVR.Range("N62").Value = -257871.25

# WVR row 65
$ws_WVR.Range("H65").Value = 172748
$ws_WVR.Range("I65").Value = 4997.5
$ws_WVR.Range("J65").Value = 256623.25
$ws_WVR.Range("K65").Value = 24987.5
$ws_WVR.Range("L65").Value = 1283116.25
$ws_WVR.Range("M65").Value = -21867.5
$ws_WVR.Range("N65").Value = -1289356.25

# WVR row 122
$ws_WVR.Range("H122").Value = 12881.223
$ws_WVR.Range("I122").Value = 17819.334
$ws_WVR.Range("J122").Value = 3005
$ws_WVR.Range("K122").Value = 53458.00199999999
$ws_WVR.Range("L122").Value = 9015
$ws_WVR.Range("M122").Value = -51008.00199999999
$ws_WVR.Range("N122").Value = -13915

# WVR row 126
$ws_WVR.Range("H126").Value = 54501.25
$ws_WVR.Range("I126").Value = 75866.07000000001
$ws_WVR.Range("J126").Value = 4650
$ws_WVR.Range("K126").Value = 227598.21
$ws_WVR.Range("L126").Value = 13950
$ws_WVR.Range("M126").Value = -225128.21
$ws_WVR.Range("N126").Value = -18890

# WVR row 132
$ws_WVR.Range("H132").Value = 59020.61
$ws_WVR.Range("I132").Value = 66148.25
$ws_WVR.Range("J132").Value = 1999.5
$ws_WVR.Range("K132").Value = 198444.75
$ws_WVR.Range("L132").Value = 5998.5
$ws_WVR.Range("M132").Value = -195914.75
$ws_WVR.Range("N132").Value = -11058.5

# WVR row 133
$ws_WVR.Range("H133").Value = 0
$ws_WVR.Range("J133").Value = 0
$ws_WVR.Range("L133").Value = 0
$ws_WVR.Range("N133").ClearContents()

# WVR row 136
$ws_WVR.Range("H136").Value = 1494.7646
$ws_WVR.Range("J136").Value = 2964
$ws_WVR.Range("L136").Value = 8892
$ws_WVR.Range("N136").Value = -13992

